$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.111.77'
$ws.Range("E2").Value = '  -2.22%  '
$ws.Range("D3").Value = '3.522.40'
$ws.Range("E3").Value = '  -3.29%  '
$ws.Range("E4").Value = '  -0.06%  '
$__s = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '588.74'
$ws.Range("D5").Style = $__s
$ws.Range("E5").Value = '  +1.23%  '
$__s = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '171.73'
$ws.Range("D6").Style = $__s
$ws.Range("E6").Value = '  -2.20%  '
$__s = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.615'
$ws.Range("D7").Style = $__s
$ws.Range("E7").Value = '  +0.92%  '
$ws.Range("D8").Value = '3.517.74'
$ws.Range("E8").Value = '  -3.22%  '
$ws.Range("E9").Value = '  -0.01%  '
$__s = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.189'
$ws.Range("D10").Style = $__s
$ws.Range("E10").Value = '  -3.72%  '
$__s = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.89'
$ws.Range("D11").Style = $__s
$ws.Range("E11").Value = '  -0.67%  '
$__s = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.582'
$ws.Range("D12").Style = $__s
$ws.Range("E12").Value = '  -4.28%  '
$__s = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '47.32'
$ws.Range("D13").Style = $__s
$ws.Range("E13").Value = '  -2.37%  '
$__s = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000277'
$ws.Range("D14").Style = $__s
$ws.Range("E14").Value = '  -2.75%  '
$ws.Range("D15").Value = '4.078.82'
$ws.Range("E15").Value = '  -3.55%  '
$ws.Range("B16").Value = 'Polkadot'
$ws.Range("C16").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$__s = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '8.49'
$ws.Range("D16").Style = $__s
$ws.Range("E16").Value = '  -4.67%  '
$ws.Range("B17").Value = 'BitcoinCash'
$ws.Range("C17").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$__s = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '623.97'
$ws.Range("D17").Style = $__s
$ws.Range("E17").Value = '  -6.46%  '
$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").Value = '69.162.90'
$ws.Range("E18").Value = '  -2.26%  '
$ws.Range("B19").Value = 'WrappedEther'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D19").Value = '3.513.63'
$ws.Range("E19").Value = '  -3.45%  '
$ws.Range("E20").Value = '  -0.09%  '
$__s = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.43'
$ws.Range("D21").Style = $__s
$ws.Range("E21").Value = '  -2.00%  '
$__s = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '11.17'
$ws.Range("D22").Style = $__s
$ws.Range("E22").Value = '  -2.25%  '
$__s = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.888'
$ws.Range("D23").Style = $__s
$ws.Range("E23").Value = '  -5.61%  '
$__s = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '15.97'
$ws.Range("D24").Style = $__s
$ws.Range("E24").Value = '  -6.63%  '
$__s = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '97.10'
$ws.Range("D25").Style = $__s
$ws.Range("E25").Value = '  -2.64%  '
$__s = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.84'
$ws.Range("D26").Style = $__s
$ws.Range("E26").Value = '  -2.00%  '
$ws.Range("E27").Value = '  +0.02%  '
$__s = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.64'
$ws.Range("D28").Style = $__s
$ws.Range("E28").Value = '  -5.44%  '
$__s = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.30'
$ws.Range("D29").Style = $__s
$ws.Range("E29").Value = '  -6.42%  '
$__s = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '32.82'
$ws.Range("D30").Style = $__s
$ws.Range("E30").Value = '  -5.36%  '
$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$__s = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.56'
$ws.Range("D31").Style = $__s
$ws.Range("E31").Value = '  -4.79%  '
$ws.Range("B32").Value = 'Stacks'
$ws.Range("C32").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$__s = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.14'
$ws.Range("D32").Style = $__s
$ws.Range("E32").Value = '  -5.83%  '
$ws.Range("E33").Value = '  -5.14%  '
$__s = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.96'
$ws.Range("D34").Style = $__s
$ws.Range("E34").Value = '  -7.16%  '
$__s = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '638.54'
$ws.Range("D35").Style = $__s
$ws.Range("E35").Value = '  +9.05%  '
$__s = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '10.79'
$ws.Range("D36").Style = $__s
$ws.Range("E36").Value = '  -2.54%  '
$__s = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.51'
$ws.Range("D37").Style = $__s
$ws.Range("E37").Value = '  -12.73%  '
$__s = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.103'
$ws.Range("D38").Style = $__s
$ws.Range("E38").Value = '  -4.19%  '
$__s = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '57.23'
$ws.Range("D39").Style = $__s
$ws.Range("E39").Value = '  -1.82%  '
$__s = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.999'
$ws.Range("D40").Style = $__s
$ws.Range("E40").Value = '  -0.05%  '
$__s = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0455'
$ws.Range("D41").Style = $__s
$ws.Range("E41").Value = '  -0.26%  '
$__s = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.136'
$ws.Range("D42").Style = $__s
$ws.Range("E42").Value = '  -3.62%  '
$ws.Range("D43").Value = '3.386.74'
$ws.Range("E43").Value = '  -5.08%  '
$__s = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.329'
$ws.Range("D44").Style = $__s
$ws.Range("E44").Value = '  -4.68%  '
$__s = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '32.90'
$ws.Range("D45").Style = $__s
$ws.Range("E45").Value = '  -4.86%  '
$ws.Range("D46").Value = '0.0₃0701'
$ws.Range("E46").Value = '  -4.64%  '
$__s = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.55'
$ws.Range("D47").Style = $__s
$ws.Range("E47").Value = '  -5.36%  '
$__s = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.76'
$ws.Range("D48").Style = $__s
$ws.Range("E48").Value = '  -5.52%  '
$ws.Range("E49").Value = '  -2.41%  '
$__s = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '133.51'
$ws.Range("D50").Style = $__s
$ws.Range("E50").Value = '  -1.81%  '
$__s = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '5.62'
$ws.Range("D51").Style = $__s
$ws.Range("E51").Value = '  +13.57%  '
